# Generate Report for Handoff
# Adds three newly-discovered files (one .md + two dependent .png files) to
# the localization status report: the Overview sheet and the per-language
# (zh-cn / de-de) detail sheets, each gaining rows 5-7 with status
# "Ready for handoff".

$wb = $excel.ActiveWorkbook

$srcRepo   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58c5a73005f7910d7bdb10cb9ba495aa423c5a43/e2e"
$zhcnRepo  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/aff768ab831679e26aedaf604cbff65f8ca0ba62/e2e"
$dedeRepo  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f6abcedc5de2acecebbc60be41d1a47e3ba062f0/e2e"

# ---------------------------------------------------------------------
# New file identities
# ---------------------------------------------------------------------
$mdName   = "62194a36-2259-448a-ba2f-fd664eefcda6.md"
$mdPath   = "e2e\62194a36-2259-448a-ba2f-fd664eefcda6.md"

$png1Name = "f461a675-26b1-45b1-b00d-1eb4b89334f4.png"
$png1Path = "e2e\f461a675-26b1-45b1-b00d-1eb4b89334f4.png"

$png2Name = "fb564dca-b56b-4291-97f8-f5f55a6de0df.png"
$png2Path = "e2e\fb564dca-b56b-4291-97f8-f5f55a6de0df.png"

$status     = "Ready for handoff"
$hoDate     = "2016-10-14 09:15:59"
$xliffDateZh = "2016-10-14 09:15:49"
$epoch      = "0001-01-01 00:00:00"

$zhXlf = "62194a36-2259-448a-ba2f-fd664eefcda6.c7ffa70f2ad1bf63fbecf7cdda89a9c5d70e0faf.zh-cn.xlf"
$deXlf = "62194a36-2259-448a-ba2f-fd664eefcda6.c7ffa70f2ad1bf63fbecf7cdda89a9c5d70e0faf.de-de.xlf"
$png1Xlf = "08023e7e19be37706f500c5d91b7a8f4e098cec7.png"
$png2Xlf = "bb086964243becc3b169b9fff546eba93c641518.png"

# =======================================================================
# Overview sheet: add rows 5-7
# =======================================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G7"))

$wsOverview.Range("A5").Value = $mdName
$wsOverview.Range("B5").Value = $mdPath
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = $status
$wsOverview.Range("F5").Value = $status
$wsOverview.Range("G5").Value = $hoDate

$wsOverview.Range("A6").Value = $png1Name
$wsOverview.Range("B6").Value = $png1Path
$wsOverview.Range("C6").Value = ".png"
$wsOverview.Range("D6").Value = ""
$wsOverview.Range("E6").Value = $status
$wsOverview.Range("F6").Value = $status
$wsOverview.Range("G6").Value = $hoDate

$wsOverview.Range("A7").Value = $png2Name
$wsOverview.Range("B7").Value = $png2Path
$wsOverview.Range("C7").Value = ".png"
$wsOverview.Range("D7").Value = ""
$wsOverview.Range("E7").Value = $status
$wsOverview.Range("F7").Value = $status
$wsOverview.Range("G7").Value = $hoDate

# Style + hyperlink for the "Path And Name" column (matches rows 2-4)
$wsOverview.Range("B5").Style = "Hyperlink"
$wsOverview.Range("B6").Style = "Hyperlink"
$wsOverview.Range("B7").Style = "Hyperlink"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "$srcRepo/$mdName", "", "", $mdPath) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "$srcRepo/$png1Name", "", "", $png1Path) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "$srcRepo/$png2Name", "", "", $png2Path) | Out-Null

# =======================================================================
# zh-cn sheet: add rows 5-7
# =======================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P7"))

# Row 5: the .md file itself
$wsZh.Range("A5").Value = $mdName
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = $status
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = $zhXlf
$wsZh.Range("H5").Value = $xliffDateZh
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = $epoch
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "False"
$wsZh.Range("P5").Value = ""

# Row 6: first dependent .png
$wsZh.Range("A6").Value = $png1Name
$wsZh.Range("B6").Value = ".png"
$wsZh.Range("C6").Value = $status
$wsZh.Range("D6").Value = "e2e"
$wsZh.Range("E6").Value = "ht"
$wsZh.Range("F6").Value = "False"
$wsZh.Range("G6").Value = $png1Xlf
$wsZh.Range("H6").Value = $xliffDateZh
$wsZh.Range("I6").Value = ""
$wsZh.Range("J6").Value = ""
$wsZh.Range("K6").Value = $epoch
$wsZh.Range("L6").Value = ""
$wsZh.Range("M6").Value = "True(Dependency)"
$wsZh.Range("N6").Value = $mdPath
$wsZh.Range("O6").Value = "False"
$wsZh.Range("P6").Value = ""

# Row 7: second dependent .png
$wsZh.Range("A7").Value = $png2Name
$wsZh.Range("B7").Value = ".png"
$wsZh.Range("C7").Value = $status
$wsZh.Range("D7").Value = "e2e"
$wsZh.Range("E7").Value = "ht"
$wsZh.Range("F7").Value = "False"
$wsZh.Range("G7").Value = $png2Xlf
$wsZh.Range("H7").Value = $xliffDateZh
$wsZh.Range("I7").Value = ""
$wsZh.Range("J7").Value = ""
$wsZh.Range("K7").Value = $epoch
$wsZh.Range("L7").Value = ""
$wsZh.Range("M7").Value = "True(Dependency)"
$wsZh.Range("N7").Value = $mdPath
$wsZh.Range("O7").Value = "False"
$wsZh.Range("P7").Value = ""

$wsZh.Range("A5").Style = "Hyperlink"
$wsZh.Range("A6").Style = "Hyperlink"
$wsZh.Range("A7").Style = "Hyperlink"
$wsZh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "$srcRepo/$mdName", "", "", $mdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "$srcRepo/$png1Name", "", "", $png1Name) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "$srcRepo/$png2Name", "", "", $png2Name) | Out-Null

$wsZh.Columns.Item(13).ColumnWidth = 16.85
$wsZh.Columns.Item(14).ColumnWidth = 39.16666666666667

# =======================================================================
# de-de sheet: add rows 5-7
# =======================================================================
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P7"))

# Row 5: the .md file itself
$wsDe.Range("A5").Value = $mdName
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = $status
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = $deXlf
$wsDe.Range("H5").Value = $hoDate
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = $epoch
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "False"
$wsDe.Range("P5").Value = ""

# Row 6: first dependent .png
$wsDe.Range("A6").Value = $png1Name
$wsDe.Range("B6").Value = ".png"
$wsDe.Range("C6").Value = $status
$wsDe.Range("D6").Value = "e2e"
$wsDe.Range("E6").Value = "ht"
$wsDe.Range("F6").Value = "False"
$wsDe.Range("G6").Value = $png1Xlf
$wsDe.Range("H6").Value = $hoDate
$wsDe.Range("I6").Value = ""
$wsDe.Range("J6").Value = ""
$wsDe.Range("K6").Value = $epoch
$wsDe.Range("L6").Value = ""
$wsDe.Range("M6").Value = "True(Dependency)"
$wsDe.Range("N6").Value = $mdPath
$wsDe.Range("O6").Value = "False"
$wsDe.Range("P6").Value = ""

# Row 7: second dependent .png
$wsDe.Range("A7").Value = $png2Name
$wsDe.Range("B7").Value = ".png"
$wsDe.Range("C7").Value = $status
$wsDe.Range("D7").Value = "e2e"
$wsDe.Range("E7").Value = "ht"
$wsDe.Range("F7").Value = "False"
$wsDe.Range("G7").Value = $png2Xlf
$wsDe.Range("H7").Value = $hoDate
$wsDe.Range("I7").Value = ""
$wsDe.Range("J7").Value = ""
$wsDe.Range("K7").Value = $epoch
$wsDe.Range("L7").Value = ""
$wsDe.Range("M7").Value = "True(Dependency)"
$wsDe.Range("N7").Value = $mdPath
$wsDe.Range("O7").Value = "False"
$wsDe.Range("P7").Value = ""

$wsDe.Range("A5").Style = "Hyperlink"
$wsDe.Range("A6").Style = "Hyperlink"
$wsDe.Range("A7").Style = "Hyperlink"
$wsDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "$srcRepo/$mdName", "", "", $mdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "$srcRepo/$png1Name", "", "", $png1Name) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "$srcRepo/$png2Name", "", "", $png2Name) | Out-Null

$wsDe.Columns.Item(13).ColumnWidth = 16.85
$wsDe.Columns.Item(14).ColumnWidth = 39.16666666666667
